$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8, shifting existing rows 8-26 down to 9-27
$ws.Rows.Item(8).Insert()

# Fill in the new row 8 with data (same fixed columns as the rest of the sheet,
# plus the new record's specific values)
$ws.Cells.Item(8, 1).Value = 10
$ws.Cells.Item(8, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(8, 3).Value = "La Araucanía"
$ws.Cells.Item(8, 4).Value = 44526
$ws.Cells.Item(8, 5).Value = 9
$ws.Cells.Item(8, 6).Value = "Fruta"
$ws.Cells.Item(8, 7).Value = 100103
$ws.Cells.Item(8, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(8, 9).Value = 100103003
$ws.Cells.Item(8, 10).Value = "Damasco"
$ws.Cells.Item(8, 11).Value = "Castle Brite"
$ws.Cells.Item(8, 12).Value = "Primera"
$ws.Cells.Item(8, 13).Value = 65
$ws.Cells.Item(8, 14).Value = 20000
$ws.Cells.Item(8, 15).Value = 20000
$ws.Cells.Item(8, 16).Value = 20000
$ws.Cells.Item(8, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(8, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(8, 19).Value = 2000
$ws.Cells.Item(8, 20).Value = 10
